# Updates the cryptos list (Coin / Link / Price / Volume(1h)) to reflect
# the latest scrape. Price cells (column D) are written with a guarded
# Text NumberFormat + ClearFormats() sequence so numeric-looking prices
# (e.g. "138.97") stay stored as plain text, matching how the sheet
# already represents every price as a string (not a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.278.18"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").Value = "2.382.98"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.97"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.43%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +0.96%  "

$ws.Range("D9").Value = "2.383.64"
$ws.Range("E9").Value = "  +0.70%  "

$ws.Range("E10").Value = "  -1.84%  "

$ws.Range("E11").Value = "  -1.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.12"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.65%  "

$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.830.90"
$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000168"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.19%  "

$ws.Range("D17").Value = "60.129.56"
$ws.Range("E17").Value = "  -0.73%  "

$ws.Range("D18").Value = "2.387.77"
$ws.Range("E18").Value = "  +0.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.08"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +11.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.55"
$ws.Range("D20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.80"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("E22").Value = "  +0.97%  "

$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("E25").Value = "  -1.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.28"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "562.70"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.11"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.61%  "

$ws.Range("D30").Value = "0.0₃0928"
$ws.Range("E30").Value = "  +1.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.01"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.09%  "

$ws.Range("E32").Value = "  -1.81%  "

$ws.Range("E33").Value = "  -1.69%  "

$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.71"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.61%  "

$ws.Range("E38").Value = "  -0.55%  "

$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("E40").Value = "  +0.39%  "

$ws.Range("E41").Value = "  +0.17%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.67"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.97%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.63"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.43"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.35%  "

$ws.Range("D46").Value = "0.0₆0291"
$ws.Range("E46").Value = "  +4.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.47"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.54"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.09%  "

$ws.Range("E49").Value = "  +0.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.20"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.47%  "

